$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 341.0909
$ws.Range("I12").Value = 140.5
$ws.Range("J12").Value = 455.7143
$ws.Range("K12").Value = 140.5
$ws.Range("L12").Value = 455.7143
$ws.Range("M12").Value = 29.5
$ws.Range("N12").Value = -795.7143

$ws.Range("H86").Value = 10677.23
$ws.Range("I86").Value = 7036.727
$ws.Range("J86").Value = 15388.471
$ws.Range("K86").Value = 7036.727
$ws.Range("L86").Value = 15388.471
$ws.Range("M86").Value = -5913.727
$ws.Range("N86").Value = -17634.471

$ws.Range("H89").Value = 10677.23
$ws.Range("I89").Value = 7036.727
$ws.Range("J89").Value = 15388.471
$ws.Range("K89").Value = 35183.635
$ws.Range("L89").Value = 76942.355
$ws.Range("M89").Value = -29567.635
$ws.Range("N89").Value = -88174.355

$ws.Range("H92").Value = 2924246.8
$ws.Range("I92").Value = 3472405
$ws.Range("J92").Value = 735
$ws.Range("K92").Value = 3472405
$ws.Range("L92").Value = 735
$ws.Range("M92").Value = -3471157
$ws.Range("N92").Value = -3231

$ws.Range("H132").Value = 12827587
$ws.Range("I132").Value = 13165129
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 39495387
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -39492857
$ws.Range("N132").Value = -8060

$ws.Range("H137").Value = 1728.5476
$ws.Range("I137").Value = 1758
$ws.Range("J137").Value = 1685.2354
$ws.Range("K137").Value = 5274
$ws.Range("L137").Value = 5055.706200000001
$ws.Range("M137").Value = -2724
$ws.Range("N137").Value = -10155.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4835714
$ws.Range("I32").Value = 3974.7
$ws.Range("J32").Value = 37047310
$ws.Range("K32").Value = 3974.7
$ws.Range("L32").Value = 37047310
$ws.Range("M32").Value = -3687.7
$ws.Range("N32").Value = -37047884

$ws.Range("H61").Value = 70314740
$ws.Range("I61").Value = 102274570
$ws.Range("J61").Value = 3094
$ws.Range("K61").Value = 102274570
$ws.Range("L61").Value = 3094
$ws.Range("M61").Value = -102274358
$ws.Range("N61").Value = -3518

$ws.Range("H132").Value = 1732531.6
$ws.Range("I132").Value = 2310.2222
$ws.Range("J132").Value = 3679030.5
$ws.Range("K132").Value = 6930.6666
$ws.Range("L132").Value = 11037091.5
$ws.Range("M132").Value = -4400.6666
$ws.Range("N132").Value = -11042151.5

$ws.Range("H136").Value = 70314740
$ws.Range("I136").Value = 102274570
$ws.Range("J136").Value = 3094
$ws.Range("K136").Value = 306823710
$ws.Range("L136").Value = 9282
$ws.Range("M136").Value = -306821160
$ws.Range("N136").Value = -14382

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 807.0417
$ws.Range("I94").Value = 798.3684
$ws.Range("J94").Value = 840
$ws.Range("K94").Value = 798.3684
$ws.Range("L94").Value = 840
$ws.Range("M94").Value = -347.3684
$ws.Range("N94").Value = -1742

$ws.Range("H132").Value = 50535.8
$ws.Range("J132").Value = 50535.8
$ws.Range("L132").Value = 50535.8
$ws.Range("N132").Value = -60655.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 83335110
$ws.Range("I99").Value = 111112770
$ws.Range("J99").Value = 2133.3333
$ws.Range("K99").Value = 111112770
$ws.Range("L99").Value = 2133.3333
$ws.Range("M99").Value = -111111272
$ws.Range("N99").Value = -5129.3333

$ws.Range("H103").Value = 7712
$ws.Range("I103").Value = 3424
$ws.Range("J103").Value = 12000
$ws.Range("K103").Value = 3424
$ws.Range("L103").Value = 12000
$ws.Range("M103").Value = -2252
$ws.Range("N103").Value = -14344

$ws.Range("H126").Value = 83335110
$ws.Range("I126").Value = 111112770
$ws.Range("J126").Value = 2133.3333
$ws.Range("K126").Value = 333338310
$ws.Range("L126").Value = 6399.999899999999
$ws.Range("M126").Value = -333335840
$ws.Range("N126").Value = -11339.9999

$ws.Range("H134").Value = 38462940
$ws.Range("I134").Value = 1222.4
$ws.Range("J134").Value = 166668670
$ws.Range("K134").Value = 3667.2
$ws.Range("L134").Value = 500006010
$ws.Range("M134").Value = -1132.2
$ws.Range("N134").Value = -500011080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 145.25
$ws.Range("I6").Value = 145.25
$ws.Range("K6").Value = 435.75
$ws.Range("M6").Value = -322.75

$ws.Range("H132").Value = 5319.6665
$ws.Range("I132").Value = 439.33334
$ws.Range("J132").Value = 10200
$ws.Range("K132").Value = 3954.00006
$ws.Range("L132").Value = 91800
$ws.Range("M132").Value = -1424.00006
$ws.Range("N132").Value = -96860

$ws.Range("H133").Value = 1230
$ws.Range("I133").Value = 1230
$ws.Range("K133").Value = 3690
$ws.Range("M133").Value = 1370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2300.3333
$ws.Range("I31").Value = 1450.5
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1450.5
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1158.5
$ws.Range("N31").Value = -4584

$ws.Range("H37").Value = 2300.3333
$ws.Range("I37").Value = 1450.5
$ws.Range("J37").Value = 4000
$ws.Range("K37").Value = 1450.5
$ws.Range("L37").Value = 4000
$ws.Range("M37").Value = -1173.5
$ws.Range("N37").Value = -4554

$ws.Range("H64").Value = 19789.5
$ws.Range("J64").Value = 19789.5
$ws.Range("L64").Value = 19789.5
$ws.Range("N64").Value = -20285.5

$ws.Range("H67").Value = 19789.5
$ws.Range("J67").Value = 19789.5
$ws.Range("L67").Value = 19789.5
$ws.Range("N67").Value = -21505.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 59799.5
$ws.Range("J6").Value = 59799.5
$ws.Range("L6").Value = 59799.5
$ws.Range("N6").Value = -60023.5

$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 2500
$ws.Range("K7").Value = 2500
$ws.Range("M7").Value = -2388

$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 73429.14
$ws.Range("I122").Value = 100800.8
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 302402.4
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -299952.4
$ws.Range("N122").Value = -19900

Write-Output "Applied Gungnir_Profits updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
